$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Product ID"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "List Price"
$ws.Range("D1").Value = "Currency"
$ws.Range("E1").Value = "Vendor"

$ws.Range("B15").Select()
